$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last four KPI header labels (E1:H1) were stored with a stray leading
# space in the shared strings table; rewrite them without it.
$ws.Range("E1").Value = "PassRateFirstInterview"
$ws.Range("F1").Value = "CostPerHire"
$ws.Range("G1").Value = "WorkPerformanceScore"
$ws.Range("H1").Value = "ConsciousnessScore"

# EmployeeId (column B, rows 2-11) used placeholder ids 1001-1010;
# replace them with simple sequential ids 1-10.
$employeeIds = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
for ($i = 0; $i -lt $employeeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $employeeIds[$i]
}

# Leave the selection where the author left it before saving.
$ws.Range("N16").Select()
